# chore: update Sheets via scheduled runner
# Applies refreshed market-price snapshot values (currentAveragePrice* / LevePrice* /
# LeveProfit* columns) to the affected Leve rows across the eight job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6579.7
$ws.Range("I51").Value = 5160
$ws.Range("K51").Value = 5160
$ws.Range("M51").Value = -4676

$ws.Range("H69").Value = 33564.707
$ws.Range("I69").Value = 64400
$ws.Range("J69").Value = 24076.924
$ws.Range("K69").Value = 193200
$ws.Range("L69").Value = 72230.772
$ws.Range("M69").Value = -192326
$ws.Range("N69").Value = -73978.772

$ws.Range("H72").Value = 33564.707
$ws.Range("I72").Value = 64400
$ws.Range("J72").Value = 24076.924
$ws.Range("K72").Value = 579600
$ws.Range("L72").Value = 216692.316
$ws.Range("M72").Value = -575232
$ws.Range("N72").Value = -225428.316

$ws.Range("H107").Value = 447.6
$ws.Range("I107").Value = 79.666664
$ws.Range("K107").Value = 79.666664
$ws.Range("M107").Value = 1840.333336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1700
$ws.Range("I21").Value = 1700
$ws.Range("K21").Value = 1700
$ws.Range("M21").Value = -1326

$ws.Range("H32").Value = 1910
$ws.Range("I32").Value = 1910
$ws.Range("K32").Value = 1910
$ws.Range("M32").Value = -1623

$ws.Range("H45").Value = 4246.4
$ws.Range("I45").Value = 4246.4
$ws.Range("K45").Value = 4246.4
$ws.Range("M45").Value = -3869.4

$ws.Range("H63").Value = 11459.883
$ws.Range("I63").Value = 12582.1
$ws.Range("K63").Value = 12582.1
$ws.Range("M63").Value = -11896.1

$ws.Range("H66").Value = 11459.883
$ws.Range("I66").Value = 12582.1
$ws.Range("K66").Value = 62910.5
$ws.Range("M66").Value = -59478.5

$ws.Range("H110").Value = 3732.818
$ws.Range("I110").Value = 4140.222
$ws.Range("K110").Value = 4140.222
$ws.Range("M110").Value = -2095.222

$ws.Range("H122").Value = 1639.2222
$ws.Range("I122").Value = 882
$ws.Range("K122").Value = 2646
$ws.Range("M122").Value = -196

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H20").Value = 1501.5834
$ws.Range("J20").Value = 1152.25
$ws.Range("L20").Value = 1152.25
$ws.Range("N20").Value = -1646.25

$ws.Range("H86").Value = 1773.4
$ws.Range("J86").Value = 1756
$ws.Range("L86").Value = 1756
$ws.Range("N86").Value = -4002

$ws.Range("H89").Value = 1773.4
$ws.Range("J89").Value = 1756
$ws.Range("L89").Value = 8780
$ws.Range("N89").Value = -20012

$ws.Range("H105").Value = 2880.25
$ws.Range("J105").Value = 2823.5
$ws.Range("L105").Value = 2823.5
$ws.Range("N105").Value = -6317.5

$ws.Range("H107").Value = 5000
$ws.Range("I107").Value = 5000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 5000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -3080
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3103.7778
$ws.Range("I16").Value = 2787.8
$ws.Range("J16").Value = 3498.75
$ws.Range("K16").Value = 2787.8
$ws.Range("L16").Value = 3498.75
$ws.Range("M16").Value = -2500.8
$ws.Range("N16").Value = -4072.75

$ws.Range("H31").Value = 4285.4287
$ws.Range("J31").Value = 1999.75
$ws.Range("L31").Value = 1999.75
$ws.Range("N31").Value = -2589.75

$ws.Range("H34").Value = 4285.4287
$ws.Range("J34").Value = 1999.75
$ws.Range("L34").Value = 1999.75
$ws.Range("N34").Value = -2403.75

$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

$ws.Range("H105").Value = 2851.1667
$ws.Range("I105").Value = 2026.75
$ws.Range("K105").Value = 2026.75
$ws.Range("M105").Value = -279.75

$ws.Range("H113").Value = 3103.7778
$ws.Range("I113").Value = 2787.8
$ws.Range("J113").Value = 3498.75
$ws.Range("K113").Value = 2787.8
$ws.Range("L113").Value = 3498.75
$ws.Range("M113").Value = -617.8000000000002
$ws.Range("N113").Value = -7838.75

$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws.Range("H132").Value = 3310.25
$ws.Range("I132").Value = 3643.3333
$ws.Range("K132").Value = 10929.9999
$ws.Range("M132").Value = -8399.999899999999

$ws.Range("H133").Value = 38442.6
$ws.Range("I133").Value = 25296
$ws.Range("K133").Value = 25296
$ws.Range("M133").Value = -22766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 960.8889
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2615.1667
$ws.Range("I80").Value = 2538.2
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2538.2
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1540.2
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 2615.1667
$ws.Range("I83").Value = 2538.2
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 12691
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -7699
$ws.Range("N83").Value = -24984

$ws.Range("H102").Value = 5108.2856
$ws.Range("I102").Value = 4951.6
$ws.Range("K102").Value = 4951.6
$ws.Range("M102").Value = -3329.6

$ws.Range("H113").Value = 3933.3333
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 5400
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 5400
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -9740

$ws.Range("H126").Value = 6313.1665
$ws.Range("I126").Value = 6275.8
$ws.Range("K126").Value = 18827.4
$ws.Range("M126").Value = -16357.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9075.385
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 9075.385
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 9075.385
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -9299.385

$ws.Range("H46").Value = 85850.164
$ws.Range("I46").Value = 126525.25
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 126525.25
$ws.Range("L46").Value = 4500
$ws.Range("M46").Value = -126337.25
$ws.Range("N46").Value = -4876

$ws.Range("H122").Value = 10699.2
$ws.Range("I122").Value = 14000.8
$ws.Range("K122").Value = 42002.39999999999
$ws.Range("M122").Value = -39552.39999999999

$ws.Range("H126").Value = 9075.385
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 9075.385
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 27226.155
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -32166.155

$ws.Range("H132").Value = 3596.4443
$ws.Range("I132").Value = 3786.5908
$ws.Range("K132").Value = 11359.7724
$ws.Range("M132").Value = -8829.7724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 44999
$ws.Range("J49").Value = 44999
$ws.Range("L49").Value = 44999
$ws.Range("N49").Value = -45459

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H107").Value = 915
$ws.Range("I107").Value = 889.8
$ws.Range("K107").Value = 2669.4
$ws.Range("M107").Value = -749.3999999999996

$ws.Range("H113").Value = 951.1
$ws.Range("I113").Value = 939.7143
$ws.Range("K113").Value = 2819.1429
$ws.Range("M113").Value = -649.1428999999998

$ws.Range("H136").Value = 7432.2666
$ws.Range("I136").Value = 8394.583000000001
$ws.Range("K136").Value = 25183.749
$ws.Range("M136").Value = -22633.749
